# Presentation Task 11, Scrum Sprint 2
# Applies the sprint-2 backlog updates: new "overview / display" stories in
# the Sprint Backlog, actual-effort figures + status updates for the
# existing sprint-1 stories, and refreshed sheet selections/active tab.

$wb = $excel.ActiveWorkbook

$wsTeam    = $wb.Worksheets.Item("ProjectTeam")
$wsBacklog = $wb.Worksheets.Item("Product Backlog")
$wsSprint  = $wb.Worksheets.Item("Sprint Backlog")

# ---------------------------------------------------------------------
# Sprint Backlog sheet: effort-actual (K) figures + status (L) updates
# for the sprint-1 stories already on the sheet.
# ---------------------------------------------------------------------
$wsSprint.Range("K2").Value = 8
$wsSprint.Range("L2").Value = "in progress"

$wsSprint.Range("K3").Value = 12

$wsSprint.Range("K4").Value = 8

$wsSprint.Range("K5").Value = 6

$wsSprint.Range("K6").Value = 0
$wsSprint.Range("L6").Value = "deferred"

# ---------------------------------------------------------------------
# New sprint-2 stories (rows 7-11). Filled column-by-column (Name,
# Components, Description, then Owner/Reviewer/Status) to match the
# order the shared strings were authored in.
# ---------------------------------------------------------------------
$wsSprint.Range("C7").Value  = "Patient display"
$wsSprint.Range("C8").Value  = "Medication display"
$wsSprint.Range("C9").Value  = "Disgnosis display"
$wsSprint.Range("C10").Value = "laboraty display"
$wsSprint.Range("C11").Value = "schedule view"

$wsSprint.Range("E7").Value  = "Database, JPA, Controller,UI"
$wsSprint.Range("E8").Value  = "Database, JPA, Controller,UI"
$wsSprint.Range("E9").Value  = "Database, JPA, Controller,UI"
$wsSprint.Range("E10").Value = "Database, JPA, Controller,UI"
$wsSprint.Range("E11").Value = "Database, JPA, Controller,UI"

$wsSprint.Range("D7").Value  = "Shows in patient overview  information about patient"
$wsSprint.Range("D8").Value  = "Shows in medication overview  information about medication"
$wsSprint.Range("D9").Value  = "Shows in disgnosis overview  information about disgnosis"
$wsSprint.Range("D10").Value = "Shows in laboraty overview  information about laboraty"
$wsSprint.Range("D11").Value = "Shows in schedule overview  information about schedule"

$wsSprint.Range("A11").Value = 2.5
$wsSprint.Range("B11").Value = 2

$wsSprint.Range("F7").Value  = "dittp2"
$wsSprint.Range("G7").Value  = "varan2"
$wsSprint.Range("H7").Value  = "high"

$wsSprint.Range("F8").Value  = "lamlr1"
$wsSprint.Range("G8").Value  = "jolop1"
$wsSprint.Range("H8").Value  = "high"

$wsSprint.Range("F9").Value  = "jolop1"
$wsSprint.Range("G9").Value  = "telec1"
$wsSprint.Range("H9").Value  = "medium"

$wsSprint.Range("F10").Value = "kammf1"
$wsSprint.Range("G10").Value = "dittp2"
$wsSprint.Range("H10").Value = "medium"

$wsSprint.Range("F11").Value = "sevib1"
$wsSprint.Range("G11").Value = "lamlr1"
$wsSprint.Range("H11").Value = "high"
$wsSprint.Range("I11").Value = 8

# Widen the Description/Components columns to fit the new text.
$wsSprint.Columns.Item(4).ColumnWidth = 56.81640625 - 0.8333333333333334
$wsSprint.Columns.Item(5).ColumnWidth = 27.26953125 - 0.8333333333333334

# ---------------------------------------------------------------------
# Sheet selections / active sheet. Selecting a range on a sheet makes it
# the active one, so the sheet that should end up active (ProjectTeam)
# is selected last.
# ---------------------------------------------------------------------
$null = $wsBacklog.Range("C17").Select()
$null = $wsSprint.Range("I10").Select()
$null = $wsTeam.Range("D17").Select()

Write-Output "Sprint 2 backlog updated"
